$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- O5:O7 data cells: copy format from column A (plain data style) and
# --- copy the literal text "FALSE" value from an existing text cell (K5)
# --- so it stays a shared string instead of being coerced to a boolean.
$ws.Range("A5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("O7").PasteSpecial(-4122)

$ws.Range("K5").Copy()
$ws.Range("O5").PasteSpecial(-4163)
$ws.Range("K6").Copy()
$ws.Range("O6").PasteSpecial(-4163)
$ws.Range("K7").Copy()
$ws.Range("O7").PasteSpecial(-4163)

# --- O4 header cell: copy format from N4 (bold header style) then bump
# --- the font size to 12pt and set the new column's label text.
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value2 = "Internal Assignment"
$ws.Range("O4").Font.Size = 12

# --- Update sheet selection to mirror the new active column
$ws.Range("O4:O7").Select()

Write-Host "done"
